$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '63.719.96'
Set-TextValue $ws.Range('E2') '  +0.72%  '
Set-TextValue $ws.Range('D3') '3.313.85'
Set-TextValue $ws.Range('E3') '  +2.58%  '
Set-TextValue $ws.Range('D5') '607.63'
Set-TextValue $ws.Range('E5') '  +2.26%  '
Set-TextValue $ws.Range('D6') '141.73'
Set-TextValue $ws.Range('E6') '  +0.42%  '
Set-TextValue $ws.Range('D8') '3.313.77'
Set-TextValue $ws.Range('E8') '  +2.68%  '
Set-TextValue $ws.Range('E9') '  +0.15%  '
Set-TextValue $ws.Range('E10') '  +1.75%  '
Set-TextValue $ws.Range('D11') '5.54'
Set-TextValue $ws.Range('E11') '  +3.44%  '
Set-TextValue $ws.Range('D12') '0.469'
Set-TextValue $ws.Range('E12') '  +1.08%  '
Set-TextValue $ws.Range('D14') '34.99'
Set-TextValue $ws.Range('E14') '  +1.87%  '
Set-TextValue $ws.Range('D15') '3.857.23'
Set-TextValue $ws.Range('E15') '  +2.57%  '
Set-TextValue $ws.Range('E16') '  +0.39%  '
Set-TextValue $ws.Range('D17') '3.313.05'
Set-TextValue $ws.Range('E17') '  +2.59%  '
Set-TextValue $ws.Range('D18') '63.783.46'
Set-TextValue $ws.Range('E18') '  +0.78%  '
Set-TextValue $ws.Range('D19') '6.87'
Set-TextValue $ws.Range('E19') '  +1.31%  '
Set-TextValue $ws.Range('D20') '481.17'
Set-TextValue $ws.Range('E20') '  +1.64%  '
Set-TextValue $ws.Range('D21') '14.08'
Set-TextValue $ws.Range('E21') '  -0.62%  '
Set-TextValue $ws.Range('D22') '0.741'
Set-TextValue $ws.Range('E22') '  +1.59%  '
Set-TextValue $ws.Range('D23') '8.00'
Set-TextValue $ws.Range('E23') '  +0.95%  '
Set-TextValue $ws.Range('E24') '  +6.40%  '
Set-TextValue $ws.Range('E25') '  +1.63%  '
Set-TextValue $ws.Range('E26') '  +0.03%  '
Set-TextValue $ws.Range('E28') '  -0.02%  '
Set-TextValue $ws.Range('D29') '8.23'
Set-TextValue $ws.Range('E29') '  +1.61%  '
Set-TextValue $ws.Range('D30') '7.17'
Set-TextValue $ws.Range('E30') '  -5.05%  '
Set-TextValue $ws.Range('E31') '  +1.69%  '
Set-TextValue $ws.Range('E32') '  +5.43%  '
Set-TextValue $ws.Range('E33') '  -0.76%  '
Set-TextValue $ws.Range('E34') '  -0.38%  '
Set-TextValue $ws.Range('D35') '1.10'
Set-TextValue $ws.Range('E35') '  +1.37%  '
Set-TextValue $ws.Range('D36') '6.07'
Set-TextValue $ws.Range('E36') '  +2.64%  '
Set-TextValue $ws.Range('D37') '52.54'
Set-TextValue $ws.Range('E37') '  -0.17%  '
Set-TextValue $ws.Range('E38') '  +6.10%  '
Set-TextValue $ws.Range('E39') '  +1.94%  '
Set-TextValue $ws.Range('B40') 'Bittensor'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D40') '434.23'
Set-TextValue $ws.Range('E40') '  +2.89%  '
Set-TextValue $ws.Range('B41') 'Maker'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D41') '3.121.52'
Set-TextValue $ws.Range('E41') '  +5.15%  '
Set-TextValue $ws.Range('E42') '  +9.66%  '
Set-TextValue $ws.Range('E43') '  -0.54%  '
Set-TextValue $ws.Range('D44') '2.74'
Set-TextValue $ws.Range('E44') '  -0.19%  '
Set-TextValue $ws.Range('E45') '  -0.15%  '
Set-TextValue $ws.Range('D46') '2.24'
Set-TextValue $ws.Range('D47') '36.84'
Set-TextValue $ws.Range('E47') '  +9.61%  '
Set-TextValue $ws.Range('D48') '26.36'
Set-TextValue $ws.Range('E48') '  +1.56%  '
Set-TextValue $ws.Range('E50') '  -1.44%  '
Set-TextValue $ws.Range('B51') 'Monero'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D51') '124.67'
Set-TextValue $ws.Range('E51') '  +3.07%  '
